$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(6,0,7,3),
    @(4,2,3,1),
    @(6,0,4,2),
    @(2,2,5,0),
    @(5,2,3,1),
    @(3,0,3,3),
    @(4,0,5,2),
    @(4,3,3,0),
    @(5,0,6,2),
    @(3,2,6,0),
    @(2,2,3,1),
    @(5,1,5,2),
    @(2,3,2,0),
    @(6,2,7,0),
    @(7,3,5,0),
    @(5,0,6,2),
    @(3,3,2,0),
    @(5,0,6,2),
    @(5,2,4,0),
    @(3,3,3,0)
)

$startRow = 1051
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
}

# Update the view to match where the new data ends (row 1071, column A).
try {
    $ws.Application.ActiveWindow.ScrollRow = 1046
} catch {}
$ws.Range("A1071").Select()
